$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.425.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.29%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.869.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.39%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.47%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7056"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.88%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07903"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3139"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.50%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.65%  "

$ws.Range("E11").Value = "  -4.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.881.92"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.37%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.201"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.93%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.80%  "

$ws.Range("E15").Value = "  -1.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.524"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.17%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008392"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.86%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.428.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "254.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.77%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.132.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.86%  "

$ws.Range("E21").Value = "  -1.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.643"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.82%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1560"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.014"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.44%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.74%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.37%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.508"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.339"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.89%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.265"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.211"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.72%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05309"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.42%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.898"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.33%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7527"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.76%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.176"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.50%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.714"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.89%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01890"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.281.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.767"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.50%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8943"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.70%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.033"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.20%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "109.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.34%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.06%  "

$ws.Range("E46").Value = "  -3.53%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.033.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.73%  "

$ws.Range("E48").Value = "  -0.21%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.589"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.27%  "

$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.5184"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.77%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4318"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.01%  "
